$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.872.41"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "1.896.32"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7920"
$ws.Range("E5").Value = "  -4.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.79"
$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3158"
$ws.Range("E8").Value = "  -3.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.35"
$ws.Range("E9").Value = "  -4.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07198"
$ws.Range("E10").Value = "  +2.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08110"
$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.590"
$ws.Range("E12").Value = "  +6.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7664"
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").Value = "1.882.57"
$ws.Range("E14").Value = "  -0.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.51"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.185"
$ws.Range("E16").Value = "  +5.96%  "

$ws.Range("D17").Value = "29.873.93"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.94"
$ws.Range("E18").Value = "  -1.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.13"
$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007788"
$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.194"
$ws.Range("E21").Value = "  +18.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").Value = "2.133.12"
$ws.Range("E23").Value = "  -0.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1667"
$ws.Range("E25").Value = "  -3.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.419"
$ws.Range("E26").Value = "  +1.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.10"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.71"
$ws.Range("E28").Value = "  -1.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.061"
$ws.Range("E29").Value = "  -1.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.403"
$ws.Range("E30").Value = "  +3.23%  "

$ws.Range("E31").Value = "  +2.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.479"
$ws.Range("E32").Value = "  +4.80%  "

$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05517"
$ws.Range("E34").Value = "  -6.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.277"
$ws.Range("E35").Value = "  +1.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7400"
$ws.Range("E36").Value = "  +1.46%  "

$ws.Range("E37").Value = "  -0.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.622"
$ws.Range("E38").Value = "  -3.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01926"
$ws.Range("E39").Value = "  +0.40%  "

$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").Value = "1.153.71"
$ws.Range("E41").Value = "  +16.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.37"
$ws.Range("E42").Value = "  +2.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4421"
$ws.Range("E43").Value = "  -0.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.888"
$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8523"
$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.41"
$ws.Range("E46").Value = "  +2.65%  "

$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.876"
$ws.Range("E48").Value = "  -1.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.966"
$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.030"
$ws.Range("E50").Value = "  +11.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.452"
$ws.Range("E51").Value = "  -1.18%  "
